$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1.0
$ws.Range("B2").Value = "CR940111"
$ws.Range("C2").Value = "Object reference not set to an instance of an object."
